$d = $word.ActiveDocument

# The final paragraph in the document currently carries the lone
# "_GoBack" bookmark. We will strip that bookmark from it (it keeps its
# Heading1 style / stays as the trailing empty heading paragraph) and then
# insert all of the new "Lessons Learned" content as fresh paragraphs
# immediately before it, re-creating the bookmark at the end of the very
# last of the newly inserted paragraphs.

$lastPara = $d.Paragraphs.Last
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$newContentXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Flyer History:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Title of window missing.  Source is missing declaration of window name in View Model</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Sort by user’s drop down box should have same visibility options as date picker box.  This makes it easier for user to understand what to do</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Database Maintainer:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Text is small, layout should be wider with bigger text to make program easier for user to navigate</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Make second and third image hidden by default, less initial information for user is easier</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Flyer Creator:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Clearing description line caused message box to display, but flyer still printed.  </w:t></w:r><w:r><w:t>Warning box should not have displayed, and flyer should not have continued if there was a warning</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Didn’t create printable flyer</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$insertionPoint.InsertXML($newContentXml)

Write-Host "Inserted lessons-learned content; paragraph count now:" $d.Paragraphs.Count
